$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B7: date changed from 43327 (2018-08-15) to 43317 (2018-08-05)
$ws.Range("B7").Value = 43317

# Row 9: new prayer entry dated 43469 (2019-01-04)
$ws.Range("B9").Value = 43469
$ws.Range("E9").Formula = "=E8-3"
$ws.Range("F9").Value = 530
$ws.Range("G9").Value = 560
$ws.Range("H9").Value = 560
$ws.Range("I9").Value = 575
$ws.Range("J9").Value = 2

# Update selected cell to C13
$ws.Range("C13").Select()
